$d = $word.ActiveDocument
$t = $d.Tables(1)

# Row 1: 99.1 -> 0M
$t.Rows(1).Cells(1).Range.Text = "0M"

# Row 2: 1.91 -> 0M
$t.Rows(2).Cells(1).Range.Text = "0M"

# Row 3: 211 -> 0M
$t.Rows(3).Cells(1).Range.Text = "0M"

# Row 4: 675 -> 712
$t.Rows(4).Cells(1).Range.Text = "712"

# Row 7: 0.01636 -> 0.01987
$t.Rows(7).Cells(1).Range.Text = "0.01987"

# Row 8: 0.00500 -> 0.00732
$t.Rows(8).Cells(1).Range.Text = "0.00732"

# Row 12: 0.45820 -> 1.90552
$t.Rows(12).Cells(1).Range.Text = "1.90552"

# Row 44: collapse multi-value tab-separated run down to "99.1"
$t.Rows(44).Cells(1).Range.Text = "99.1"

# Row 45: collapse multi-value tab-separated run down to "1.91"
$t.Rows(45).Cells(1).Range.Text = "1.91"

# Row 46: collapse multi-value tab-separated run down to "211"
$t.Rows(46).Cells(1).Range.Text = "211"
